$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second row's expected-items list no longer includes "ZARA COAT 3";
# update it to match the (now shorter) list used by the third row.
$ws.Range("A2").Value = "ADIDAS ORIGINAL,IPHONE 13 PRO"

# Reflect the last saved selection/active cell.
$ws.Range("A5").Select()
